$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.993.88'
$ws.Range("E2").Value = '  +0.16%  '
$ws.Range("D3").Value = '3.517.16'
$ws.Range("E3").Value = '  -1.30%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '591.23'
$ws.Range("E5").Value = '  -1.32%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '133.83'
$ws.Range("E6").Value = '  -1.10%  '
$ws.Range("D7").Value = '3.516.42'
$ws.Range("E7").Value = '  -1.29%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("E9").Value = '  -1.00%  '
$ws.Range("E10").Value = '  +1.49%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.19'
$ws.Range("E11").Value = '  +3.26%  '
$ws.Range("E12").Value = '  +0.08%  '
$ws.Range("D13").Value = '4.113.87'
$ws.Range("E13").Value = '  -1.39%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.60'
$ws.Range("E14").Value = '  +2.19%  '
$ws.Range("E15").Value = '  -0.35%  '
$ws.Range("E16").Value = '  +0.72%  '
$ws.Range("D17").Value = '3.515.35'
$ws.Range("E17").Value = '  -1.41%  '
$ws.Range("D18").Value = '64.987.27'
$ws.Range("E18").Value = '  +0.70%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.06'
$ws.Range("E19").Value = '  +0.42%  '
$ws.Range("E20").Value = '  -0.30%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.66'
$ws.Range("E21").Value = '  -2.68%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '391.03'
$ws.Range("E22").Value = '  +0.57%  '
$ws.Range("E23").Value = '  -0.65%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '74.72'
$ws.Range("E24").Value = '  +0.94%  '
$ws.Range("D25").Value = '3.660.06'
$ws.Range("E25").Value = '  -1.37%  '
$ws.Range("E26").Value = '  -0.06%  '
$ws.Range("E27").Value = '  -3.32%  '
$ws.Range("E28").Value = '  +8.85%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.62'
$ws.Range("E29").Value = '  -1.37%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("E30").Value = '  -0.14%  '
$ws.Range("E31").Value = '  -1.11%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.31'
$ws.Range("E32").Value = '  -1.29%  '
$ws.Range("D33").Value = '3.524.99'
$ws.Range("E33").Value = '  -1.18%  '
$ws.Range("E34").Value = '  +0.24%  '
$ws.Range("E35").Value = '  +0.02%  '
$ws.Range("E36").Value = '  +1.31%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.22'
$ws.Range("E37").Value = '  +4.71%  '
$ws.Range("E38").Value = '  +1.64%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '169.36'
$ws.Range("E39").Value = '  +0.10%  '
$ws.Range("E40").Value = '  -0.15%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0805'
$ws.Range("E41").Value = '  -0.22%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.821'
$ws.Range("E42").Value = '  -0.45%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '26.10'
$ws.Range("E43").Value = '  -3.72%  '
$ws.Range("B44").Value = 'OKB'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '42.91'
$ws.Range("E44").Value = '  +0.46%  '
$ws.Range("B45").Value = 'ONDO'
$ws.Range("C45").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.25'
$ws.Range("E45").Value = '  +3.49%  '
$ws.Range("E46").Value = '  -0.08%  '
$ws.Range("E47").Value = '  -0.62%  '
$ws.Range("E48").Value = '  +0.23%  '
$ws.Range("D49").Value = '2.454.90'
$ws.Range("E49").Value = '  -0.96%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.88'
$ws.Range("E50").Value = '  -0.27%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.898'
$ws.Range("E51").Value = '  +3.54%  '
